$d = $word.ActiveDocument
$full = $d.Content.WordOpenXML()
Write-Host "LEN before: $($full.Length)"
$modified = $full -replace '<w:contextualSpacing[^/]*/>', ''
$count = ([regex]::Matches($modified, 'contextualSpacing')).Count
Write-Host "contextualSpacing count after replace: $count"
Write-Host "LEN after: $($modified.Length)"
$d.Content.InsertXML($modified)
$full2 = $d.Content.WordOpenXML()
$count2 = ([regex]::Matches($full2, 'contextualSpacing')).Count
Write-Host "contextualSpacing count after reload: $count2"
